$wb = $excel.ActiveWorkbook

# Rename second sheet from state_milk_production -> ukvacc
$ukvacc = $wb.Worksheets.Item("state_milk_production")
$ukvacc.Name = "ukvacc"

# Add a new sheet "englandvacc" after ukvacc
$englandvacc = $wb.Worksheets.Add()
$englandvacc.Name = "englandvacc"
$englandvacc.Move($null, $wb.Worksheets.Item("ukvacc"))

# Re-fetch the sheet reference since Move() invalidates the old handle
$englandvacc = $wb.Worksheets.Item("englandvacc")

# Fill in the englandvacc sheet content
$data = @(
    @("variable", "class", "description"),
    @("areaCode", "character", "Area code of the region"),
    @("areaName", "character", "Name of the region"),
    @("areaType", "character", "region"),
    @("date", "date", "Date of data publication"),
    @("Completely Vaccinated", "double", "Total number of people fully vaccinated"),
    @("First Dose", "double", "Total number of people with one dose of vaccine"),
    @("Second Dose", "double", "Total number of people with two doses of vaccines"),
    @("pop", "double", "Population of the region (2020)"),
    @("firstperc", "double", "Percentage of population with first dose"),
    @("secondperc", "double", "Percentage of population with two doses")
)

$r = 1
foreach ($row in $data) {
    $englandvacc.Cells.Item($r, 1).Value = $row[0]
    $englandvacc.Cells.Item($r, 2).Value = $row[1]
    $englandvacc.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Style the header row (center aligned), matching the other sheets' header style
$headerStyleSource = $ukvacc.Range("A1:C1")
$headerStyleSource.Copy()
$englandvacc.Range("A1:C1").PasteSpecial(-4122)

# Autofit columns A and C to the content (column B keeps the default width)
$englandvacc.Columns.Item(1).ColumnWidth = 17.751
$englandvacc.Columns.Item(3).ColumnWidth = 38.751

# Set selection on ukvacc
$ukvacc = $wb.Worksheets.Item("ukvacc")
$ukvacc.Range("A1:C1").Select()

# Make englandvacc the active/selected sheet and set its selection
$englandvacc.Activate()
$englandvacc.Range("C3").Select()

